$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MCMCpack")
$ws.Range("D1").Value = "MCMC-based Bayesian"
$co = $ws.ChartObjects().Add(10, 10, 300, 200)
$chart = $co.Chart
$chart.ChartType = -4169
$chart.SetSourceData($ws.Range("C1:D14"))
Write-Host "added"
